# quickstart-hybrid.pptx restructuring: update cached date placeholders
# (deck was resaved two days later, so every "datetime1"/"datetimeFigureOut"
# placeholder's cached text moved from 6/17/21 -> 6/19/21) and shorten the
# two version-tag callouts from "v1.0"/"v2.0" to "v1"/"v2".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "6/17/21") {
                $tr.Text = "6/19/21"
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.Slides.Item(1).Master
Update-DatePlaceholder $master.Shapes

# Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master's date placeholder.
Update-DatePlaceholder $p.NotesMaster.Shapes

# Version tags on the single content slide: "v1.0" -> "v1", "v2.0" -> "v2".
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "v1.0") {
            $tr.Text = "v1"
        } elseif ($tr.Text -eq "v2.0") {
            $tr.Text = "v2"
        }
    }
}
